# Refresh crypto price/volume figures to match the latest upstream feed snapshot.
# Two coin pairs (Algorand/Aptos and NEARProtocol/Decentraland/Elrond) were also
# reordered in the source feed, so their Coin/Link/Price/Volume cells are rewritten
# in place to reflect the new row assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.677.80'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4920'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2931'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06725'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.895.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.17'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07257'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '90.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6742'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.017'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.678.00'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007952'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('E19').Value = '  +2.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.142.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.801'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '189.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +32.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.078'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.342'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.889'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.404'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.267'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09072'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.992'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05217'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7356'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.104'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.764'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01828'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.120'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9276'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4388'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.99'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.731'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.516'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1347'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05858'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.720'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.86%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.416'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.54%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3923'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.96%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.416'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.17%  '
